$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row before row 7, shifting the existing rows 7-38 down to 8-39
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new "coude" data point
$ws.Range("A7").Value = "coude"
$ws.Range("B7").Value = 4
$ws.Range("E7").Value = 0.19

# Update the selection / active cell to match the saved view state
$ws.Range("D6").Select()
